# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Swap the order of "Santa Lucia" / "Timor Oriental" rows (202/203) so that
#    row 202 now reads "Timor Oriental" and row 203 reads "Santa Lucia",
#    keeping each row's own numeric data untouched (both rows share identical
#    figures, so only the country labels move).
$ws.Range("A202").Value = "Timor Oriental"
$ws.Range("A203").Value = "Santa Lucia"

# 2) Update the "last updated" timestamp string in A1.
$ws.Range("A1").Value = "Datos actualizados a 8 de Agosto de 2020 a las 14:22"

# 3) Refresh per-country COVID figures (Casos totales, Nuevos casos,
#    Casos activos, Recuperados, Muertes hoy, Muertes).
$updates = @(
    @{ Row = 4;   B = 5097187; C = 1663; D = 2618025; E = 2315023; G = 45; H = 164139 }  # Estados Unidos
    @{ Row = 6;   B = 2091810; C = 4946; D = 1429325; E = 619859;  G = 48; H = 42626  }  # India
    @{ Row = 22;  B = 216346;  C = 31;   D = 197400;  E = 9692 }                          # Alemania
    @{ Row = 28;  B = 112650;  C = 267;  D = 109438;  E = 3030;   G = 2;  H = 182    }   # Catar
    @{ Row = 42;  B = 68738;   C = 124;  D = 64744;   E = 3409;   G = 2;  H = 585    }   # Bielorrusia
    @{ Row = 60;  B = 33481;   C = 105;  D = 30056;   E = 2937;   G = 5;  H = 488    }   # Azerbaiyan
    @{ Row = 62;  B = 29652;   C = 595;  D = 21006;   E = 8459;   G = 6;  H = 187    }   # Uzbekistan
    @{ Row = 79;  B = 13928;   C = 206;  D = 7706;    E = 6127;   G = 1;  H = 95     }   # Estado de Palestina
    @{ Row = 101; B = 5543;    C = 77;   D = 4817;    E = 569;    G = 2;  H = 157    }   # Croacia
    @{ Row = 133; B = 1955;    C = 3;    E = 112 }                                       # Islandia
    @{ Row = 159; B = 797;     C = 8;    E = 392 }                                       # Vietnam
    @{ Row = 205; B = 23;      C = 1;    E = 1 }                                         # Nueva Caledonia
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 2).Value = $u.B
    $ws.Cells.Item($r, 3).Value = $u.C
    if ($u.ContainsKey("D")) { $ws.Cells.Item($r, 4).Value = $u.D }
    $ws.Cells.Item($r, 5).Value = $u.E
    if ($u.ContainsKey("G")) { $ws.Cells.Item($r, 7).Value = $u.G }
    if ($u.ContainsKey("H")) { $ws.Cells.Item($r, 8).Value = $u.H }
}
